$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("D1").Value = "Country"
$ws.Range("E1").Value = "BirthMonth"
$ws.Range("F1").Value = "BirthDay"
$ws.Range("G1").Value = "BirthYear"
$ws.Range("H1").Value = "Phone"
$ws.Range("I1").Value = "UserName"
$ws.Range("J1").Value = "Email"
$ws.Range("K1").Value = "Picture"
$ws.Range("L1").Value = "Description"
$ws.Range("M1").Value = "Password"
$ws.Range("N1").Value = "ConfirmPassword"

# --- New value cells (row 2) ---
$ws.Range("D2").Value = "Bulgaria"
$ws.Range("E2").Value = "3"
$ws.Range("G2").Value = "1987"
$ws.Range("H2").Value = "0897675645"
$ws.Range("I2").Value = "lichkata456"

# --- Fix existing C2 cell: drop its numeric style, make it plain text like the others ---
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "String.Empty"

# --- Email cell + hyperlink ---
$ws.Range("J2").Value = "lichkata456@abv.bg"
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:lichkata456@abv.bg")

$ws.Range("K2").Value = "C:\Users\Iliya\Desktop\photo.jpeg"
$ws.Range("L2").Value = "ALA BALA"
$ws.Range("M2").Value = "12345678"
$ws.Range("N2").Value = "12345678"

# BirthDay entered last, out of column order (matches original authoring order)
$ws.Range("F2").Value = "23"

# --- Column sizing for the new columns (best-fit, mirroring Excel's auto column sizing) ---
$ws.Columns("H:H").ColumnWidth = 10.09
$ws.Columns("I:I").ColumnWidth = 9.92
$ws.Columns("J:J").ColumnWidth = 17.92
$ws.Columns("K:K").ColumnWidth = 31.59
$ws.Columns("L:L").ColumnWidth = 10.25
$ws.Columns("M:M").ColumnWidth = 8.59
$ws.Columns("N:N").ColumnWidth = 15.75

$ws.Range("K2").Select()
